$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Report"

# Title banner - Row 1 (blue fill across A1:S1)
$ws.Range("A1:S1").Interior.Color = 16711680
$ws.Range("B1").Value = "Report"

# Data rows 4-8 for the first (Child Tag / Text) table, columns A-C
$rows = @(
    @(0, "PUMP:HRD:100",  "Details regarding the rechargeable Lithium Polymer Battery. "),
    @(1, "PUMP:HRD:105",  "Details regarding the fuel gauge hardware for the lithium polymer battery. The battery charge shall be displayed to the user. "),
    @(2, "PUMP:HRD:1000", "Details regarding the pressure sensors for use in conjunction with the ideal gas law. "),
    @(3, "PUMP:HRD:3330", "Details regarding the size and weight of the pump. "),
    @(4, "PUMP:HRD:3350", "Details regarding the full color touchscreen. ")
)

$r = 4
foreach ($item in $rows) {
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $r = $r + 1
}

# Header for first table - Row 3, columns B:C
$ws.Range("B3").Interior.Color = 255
$ws.Range("B3").Value = "Child Tag"
$ws.Range("C3").Interior.Color = 65280
$ws.Range("C3").Value = "Text"

# Data rows 4-8 for the second (Child Tag / Parent Tag) table, columns D-F
$parentTags = @("[PUMP:HRS:100]", "[PUMP:HRS:103]", "[PUMP:HRS:1000]", "[PUMP:HRS:3330]", "[PUMP:HRS:3350]")

$r = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $item = $rows[$i]
    $ws.Cells.Item($r, 4).Value = $item[0]
    $ws.Cells.Item($r, 5).Value = $item[1]
    $ws.Cells.Item($r, 6).Value = $parentTags[$i]
    $r = $r + 1
}

# Header for second table - Row 3, columns E:F
$ws.Range("E3").Interior.Color = 255
$ws.Range("E3").Value = "Child Tag"
$ws.Range("F3").Interior.Color = 8421504
$ws.Range("F3").Value = "Parent Tag"

# Column widths (A-F)
$ws.Columns.Item(1).ColumnWidth = 2.1640625
$ws.Columns.Item(2).ColumnWidth = 14.6640625
$ws.Columns.Item(3).ColumnWidth = 106.33203125
$ws.Columns.Item(4).ColumnWidth = 2.1640625
$ws.Columns.Item(5).ColumnWidth = 14.6640625
$ws.Columns.Item(6).ColumnWidth = 15.83203125
